# Applies the cryptos-list price/volume/coin update described in the commit
# "Updated cryptos list on Sat Dec 23 19:22:20 UTC 2023 with GitHub Actions".
#
# All touched columns (B=Coin, C=Link, D=Price, E=Volume(1h)) are stored in the
# workbook as text (inlineStr), not numbers -- many Price values look numeric
# ("15.80", "0.618", "2.651.40", ...) but must keep exact text formatting
# (trailing zeros, multiple dots as thousands separators, etc). A leading
# apostrophe forces Excel to store the literal text instead of re-parsing it
# as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''43.730.51'
$ws.Range("E2").Value = '''  -0.14%  '
$ws.Range("D3").Value = '''2.296.16'
$ws.Range("E3").Value = '''  -1.01%  '
$ws.Range("E4").Value = '''  -0.06%  '
$ws.Range("D5").Value = '''101.97'
$ws.Range("E5").Value = '''  +4.75%  '
$ws.Range("D6").Value = '''270.01'
$ws.Range("E6").Value = '''  -0.32%  '
$ws.Range("D7").Value = '''0.618'
$ws.Range("E7").Value = '''  -1.12%  '
$ws.Range("E8").Value = '''  -0.07%  '
$ws.Range("E9").Value = '''  -2.39%  '
$ws.Range("D10").Value = '''45.37'
$ws.Range("E10").Value = '''  -0.84%  '
$ws.Range("D11").Value = '''0.0934'
$ws.Range("E11").Value = '''  -1.44%  '
$ws.Range("D12").Value = '''8.01'
$ws.Range("E12").Value = '''  -1.15%  '
$ws.Range("E13").Value = '''  +1.48%  '
$ws.Range("D14").Value = '''15.80'
$ws.Range("E14").Value = '''  +1.99%  '
$ws.Range("B15").Value = '''Polygon'
$ws.Range("C15").Value = '''https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '''0.858'
$ws.Range("E15").Value = '''  -1.01%  '
$ws.Range("B16").Value = '''WrappedEther'
$ws.Range("C16").Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '''2.313.19'
$ws.Range("E16").Value = '''  -0.56%  '
$ws.Range("B17").Value = '''WrappedBTC'
$ws.Range("C17").Value = '''https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '''43.745.15'
$ws.Range("E17").Value = '''  -0.12%  '
$ws.Range("B18").Value = '''ShibaInu'
$ws.Range("C18").Value = '''https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '''0.0000110'
$ws.Range("E18").Value = '''  +1.26%  '
$ws.Range("B19").Value = '''Uniswap'
$ws.Range("C19").Value = '''https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '''6.24'
$ws.Range("E19").Value = '''  -2.77%  '
$ws.Range("B20").Value = '''Litecoin'
$ws.Range("C20").Value = '''https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").Value = '''72.25'
$ws.Range("E20").Value = '''  -0.67%  '
$ws.Range("B21").Value = '''ImmutableX'
$ws.Range("C21").Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D21").Value = '''2.47'
$ws.Range("E21").Value = '''  +8.38%  '
$ws.Range("B22").Value = '''BitcoinCash'
$ws.Range("C22").Value = '''https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").Value = '''233.57'
$ws.Range("E22").Value = '''  -2.26%  '
$ws.Range("B23").Value = '''PancakeSwap'
$ws.Range("C23").Value = '''https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").Value = '''2.85'
$ws.Range("E23").Value = '''  +12.85%  '
$ws.Range("B24").Value = '''InternetComputer(DFINITY)'
$ws.Range("C24").Value = '''https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").Value = '''9.22'
$ws.Range("E24").Value = '''  -2.17%  '
$ws.Range("B25").Value = '''Dai'
$ws.Range("C25").Value = '''https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '''0.999'
$ws.Range("E25").Value = '''  -0.03%  '
$ws.Range("B26").Value = '''Cosmos'
$ws.Range("C26").Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''11.19'
$ws.Range("E26").Value = '''  -1.14%  '
$ws.Range("B27").Value = '''WEMIXToken'
$ws.Range("C27").Value = '''https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D27").Value = '''3.45'
$ws.Range("E27").Value = '''  -0.25%  '
$ws.Range("B28").Value = '''Toncoin'
$ws.Range("C28").Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '''2.29'
$ws.Range("E28").Value = '''  +0.41%  '
$ws.Range("B29").Value = '''InjectiveProtocol'
$ws.Range("C29").Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").Value = '''39.49'
$ws.Range("E29").Value = '''  +3.22%  '
$ws.Range("B30").Value = '''Monero'
$ws.Range("C30").Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = '''177.55'
$ws.Range("E30").Value = '''  +1.48%  '
$ws.Range("B31").Value = '''EthereumClassic'
$ws.Range("C31").Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '''21.80'
$ws.Range("E31").Value = '''  -2.62%  '
$ws.Range("B32").Value = '''Hedera'
$ws.Range("C32").Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.0898'
$ws.Range("E32").Value = '''  -0.65%  '
$ws.Range("B33").Value = '''Filecoin'
$ws.Range("C33").Value = '''https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''5.45'
$ws.Range("E33").Value = '''  -0.35%  '
$ws.Range("B34").Value = '''RenderToken'
$ws.Range("C34").Value = '''https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").Value = '''4.82'
$ws.Range("E34").Value = '''  +9.43%  '
$ws.Range("E35").Value = '''  +0.02%  '
$ws.Range("B36").Value = '''Kaspa'
$ws.Range("C36").Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '''0.108'
$ws.Range("E36").Value = '''  -0.08%  '
$ws.Range("B37").Value = '''VeChain'
$ws.Range("C37").Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.0353'
$ws.Range("E37").Value = '''  -1.67%  '
$ws.Range("B38").Value = '''NEARProtocol'
$ws.Range("C38").Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '''3.55'
$ws.Range("E38").Value = '''  +5.35%  '
$ws.Range("B39").Value = '''LidoDAOToken'
$ws.Range("C39").Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '''2.34'
$ws.Range("E39").Value = '''  -0.54%  '
$ws.Range("B40").Value = '''Algorand'
$ws.Range("C40").Value = '''https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '''0.234'
$ws.Range("E40").Value = '''  -3.21%  '
$ws.Range("B41").Value = '''ARBITRUM'
$ws.Range("C41").Value = '''https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '''1.39'
$ws.Range("E41").Value = '''  +1.06%  '
$ws.Range("B42").Value = '''Celestia'
$ws.Range("C42").Value = '''https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").Value = '''12.31'
$ws.Range("E42").Value = '''  +1.12%  '
$ws.Range("B43").Value = '''MultiversX'
$ws.Range("C43").Value = '''https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").Value = '''64.81'
$ws.Range("E43").Value = '''  +4.48%  '
$ws.Range("B44").Value = '''FraxShare'
$ws.Range("C44").Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '''8.84'
$ws.Range("E44").Value = '''  -3.76%  '
$ws.Range("B45").Value = '''THORChain'
$ws.Range("C45").Value = '''https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D45").Value = '''5.24'
$ws.Range("E45").Value = '''  -2.41%  '
$ws.Range("B46").Value = '''Cronos'
$ws.Range("C46").Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '''0.102'
$ws.Range("E46").Value = '''  -1.01%  '
$ws.Range("B47").Value = '''TrustWalletToken'
$ws.Range("C47").Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").Value = '''1.21'
$ws.Range("E47").Value = '''  +0.98%  '
$ws.Range("B48").Value = '''Aave'
$ws.Range("C48").Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '''98.18'
$ws.Range("E48").Value = '''  -2.19%  '
$ws.Range("B49").Value = '''WOONetwork'
$ws.Range("C49").Value = '''https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D49").Value = '''0.449'
$ws.Range("E49").Value = '''  +7.89%  '
$ws.Range("B50").Value = '''Stacks'
$ws.Range("C50").Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = '''1.53'
$ws.Range("E50").Value = '''  +12.01%  '
$ws.Range("B51").Value = '''RocketPoolETH'
$ws.Range("C51").Value = '''https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '''2.527.74'
$ws.Range("E51").Value = '''  -0.77%  '
